$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing odds values in row 2 and row 3 ---
$ws.Cells.Item(2, 7).Value = 1.3  # G2: 1.33 -> 1.3
$ws.Cells.Item(2, 8).Value = 4.5  # H2: 4.4 -> 4.5
$ws.Cells.Item(2, 9).Value = 13  # I2: 12 -> 13
$ws.Cells.Item(2, 10).Value = 1.83  # J2: 1.91 -> 1.83
$ws.Cells.Item(2, 17).Value = 2.2  # Q2: 2.3 -> 2.2
$ws.Cells.Item(2, 18).Value = 1.65  # R2: 1.6 -> 1.65
$ws.Cells.Item(2, 19).Value = 1.44  # S2: 1.5 -> 1.44
$ws.Cells.Item(2, 20).Value = 2.63  # T2: 2.5 -> 2.63
$ws.Cells.Item(2, 24).Value = 4.75  # X2: 5 -> 4.75
$ws.Cells.Item(2, 26).Value = 7  # Z2: 7.5 -> 7
$ws.Cells.Item(2, 30).Value = 9.5  # AD2: 9 -> 9.5
$ws.Cells.Item(2, 34).Value = 21  # AH2: 19 -> 21
$ws.Cells.Item(2, 39).Value = 151  # AM2: 126 -> 151
$ws.Cells.Item(2, 40).Value = 2.88  # AN2: 3 -> 2.88
$ws.Cells.Item(2, 43).Value = 19  # AQ2: 21 -> 19
$ws.Cells.Item(2, 46).Value = 2.63  # AT2: 2.5 -> 2.63
$ws.Cells.Item(2, 48).Value = 126  # AV2: 101 -> 126
$ws.Cells.Item(3, 21).Value = 1.88  # U3: 1.87 -> 1.88
$ws.Cells.Item(3, 22).Value = 1.88  # V3: 1.89 -> 1.88

# --- Append new match rows 6 and 7 ---
# Row 6
$ws.Cells.Item(6, 1).Value = "OK8mlfzI"
$ws.Cells.Item(6, 2).Value = "28/10/2024"
$ws.Cells.Item(6, 3).Value = "14:00"
$ws.Cells.Item(6, 4).Value = "TURKEY - SUPER LIG"
$ws.Cells.Item(6, 5).Value = "Galatasaray"
$ws.Cells.Item(6, 6).Value = "Besiktas"
$ws.Cells.Item(6, 7).Value = 1.8
$ws.Cells.Item(6, 8).Value = 4.2
$ws.Cells.Item(6, 9).Value = 3.8
$ws.Cells.Item(6, 10).Value = 2.25
$ws.Cells.Item(6, 11).Value = 2.6
$ws.Cells.Item(6, 12).Value = 4
$ws.Cells.Item(6, 13).Value = 1.01
$ws.Cells.Item(6, 14).Value = 23
$ws.Cells.Item(6, 15).Value = 1.11
$ws.Cells.Item(6, 16).Value = 6.5
$ws.Cells.Item(6, 17).Value = 1.4
$ws.Cells.Item(6, 18).Value = 2.88
$ws.Cells.Item(6, 19).Value = 1.22
$ws.Cells.Item(6, 20).Value = 4
$ws.Cells.Item(6, 21).Value = 1.4
$ws.Cells.Item(6, 22).Value = 2.75
$ws.Cells.Item(6, 23).Value = 13
$ws.Cells.Item(6, 24).Value = 12
$ws.Cells.Item(6, 25).Value = 9
$ws.Cells.Item(6, 26).Value = 17
$ws.Cells.Item(6, 27).Value = 12
$ws.Cells.Item(6, 28).Value = 17
$ws.Cells.Item(6, 29).Value = 23
$ws.Cells.Item(6, 30).Value = 9
$ws.Cells.Item(6, 31).Value = 12
$ws.Cells.Item(6, 32).Value = 29
$ws.Cells.Item(6, 33).Value = 81
$ws.Cells.Item(6, 34).Value = 19
$ws.Cells.Item(6, 35).Value = 26
$ws.Cells.Item(6, 36).Value = 13
$ws.Cells.Item(6, 37).Value = 41
$ws.Cells.Item(6, 38).Value = 26
$ws.Cells.Item(6, 39).Value = 23
$ws.Cells.Item(6, 40).Value = 4.5
$ws.Cells.Item(6, 41).Value = 9
$ws.Cells.Item(6, 42).Value = 15
$ws.Cells.Item(6, 43).Value = 23
$ws.Cells.Item(6, 44).Value = 34
$ws.Cells.Item(6, 45).Value = 67
$ws.Cells.Item(6, 46).Value = 4
$ws.Cells.Item(6, 47).Value = 7
$ws.Cells.Item(6, 48).Value = 34
$ws.Cells.Item(6, 49).Value = 6.5
$ws.Cells.Item(6, 50).Value = 19
$ws.Cells.Item(6, 51).Value = 21
$ws.Cells.Item(6, 52).Value = 51
$ws.Cells.Item(6, 53).Value = 51
$ws.Cells.Item(6, 54).Value = 101
$ws.Cells.Item(6, 55).Value = 251
$ws.Cells.Item(6, 56).Value = 301

# Row 7
$ws.Cells.Item(7, 1).Value = "xQiOD6C7"
$ws.Cells.Item(7, 2).Value = "28/10/2024"
$ws.Cells.Item(7, 3).Value = "14:00"
$ws.Cells.Item(7, 4).Value = "TURKEY - 1. LIG"
$ws.Cells.Item(7, 5).Value = "Manisa FK"
$ws.Cells.Item(7, 6).Value = "Erokspor"
$ws.Cells.Item(7, 7).Value = 1.75
$ws.Cells.Item(7, 8).Value = 3.7
$ws.Cells.Item(7, 9).Value = 4
$ws.Cells.Item(7, 10).Value = 2.38
$ws.Cells.Item(7, 11).Value = 2.3
$ws.Cells.Item(7, 12).Value = 4.5
$ws.Cells.Item(7, 13).Value = 1.04
$ws.Cells.Item(7, 14).Value = 13
$ws.Cells.Item(7, 15).Value = 1.22
$ws.Cells.Item(7, 16).Value = 4
$ws.Cells.Item(7, 17).Value = 1.7
$ws.Cells.Item(7, 18).Value = 2.1
$ws.Cells.Item(7, 19).Value = 1.33
$ws.Cells.Item(7, 20).Value = 3.25
$ws.Cells.Item(7, 21).Value = 1.67
$ws.Cells.Item(7, 22).Value = 2.1
$ws.Cells.Item(7, 23).Value = 8.5
$ws.Cells.Item(7, 24).Value = 9.5
$ws.Cells.Item(7, 25).Value = 8.5
$ws.Cells.Item(7, 26).Value = 15
$ws.Cells.Item(7, 27).Value = 13
$ws.Cells.Item(7, 28).Value = 23
$ws.Cells.Item(7, 29).Value = 13
$ws.Cells.Item(7, 30).Value = 7.5
$ws.Cells.Item(7, 31).Value = 15
$ws.Cells.Item(7, 32).Value = 41
$ws.Cells.Item(7, 33).Value = 151
$ws.Cells.Item(7, 34).Value = 13
$ws.Cells.Item(7, 35).Value = 23
$ws.Cells.Item(7, 36).Value = 13
$ws.Cells.Item(7, 37).Value = 41
$ws.Cells.Item(7, 38).Value = 29
$ws.Cells.Item(7, 39).Value = 34
$ws.Cells.Item(7, 40).Value = 4
$ws.Cells.Item(7, 41).Value = 9
$ws.Cells.Item(7, 42).Value = 17
$ws.Cells.Item(7, 43).Value = 29
$ws.Cells.Item(7, 44).Value = 41
$ws.Cells.Item(7, 45).Value = 101
$ws.Cells.Item(7, 46).Value = 3.25
$ws.Cells.Item(7, 47).Value = 7.5
$ws.Cells.Item(7, 48).Value = 51
$ws.Cells.Item(7, 49).Value = 6
$ws.Cells.Item(7, 50).Value = 21
$ws.Cells.Item(7, 51).Value = 26
$ws.Cells.Item(7, 52).Value = 67
$ws.Cells.Item(7, 53).Value = 81
$ws.Cells.Item(7, 54).Value = 151
$ws.Cells.Item(7, 55).Value = 126
$ws.Cells.Item(7, 56).Value = 151

